$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's row of data (row 53), continuing the daily series.
$row = 53

$ws.Cells.Item($row, 1).Value = 46002
$ws.Cells.Item($row, 2).Value = 118
$ws.Cells.Item($row, 3).Value = 128
$ws.Cells.Item($row, 4).Value = 119

# Match the date-column style used by the rest of column A (row above).
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
